$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) cells to remain plain text, matching the
# original inlineStr storage (values use "." as thousands separator
# in several rows, so Excel must not reinterpret them as numbers).
$textCells = @("D2", "D3", "D5", "D6", "D9", "D11", "D12", "D14", "D15", "D17", "D18", "D19", "D20", "D21", "D22", "D24", "D25", "D26", "D27", "D30", "D31", "D32", "D34", "D36", "D40", "D41", "D42", "D43", "D45", "D46", "D47", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values from the latest crypto data pull.
$ws.Range("D2").Value = "62.841.83"
$ws.Range("E2").Value = "  -5.24%  "
$ws.Range("D3").Value = "3.212.25"
$ws.Range("E3").Value = "  -6.38%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "174.43"
$ws.Range("E5").Value = "  -6.94%  "
$ws.Range("D6").Value = "513.05"
$ws.Range("E6").Value = "  -4.94%  "
$ws.Range("E7").Value = "  -5.01%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("D9").Value = "3.214.31"
$ws.Range("E10").Value = "  -6.95%  "
$ws.Range("D11").Value = "52.17"
$ws.Range("E11").Value = "  -11.35%  "
$ws.Range("D12").Value = "0.129"
$ws.Range("E12").Value = "  -6.43%  "
$ws.Range("E13").Value = "  -3.97%  "
$ws.Range("D14").Value = "8.85"
$ws.Range("E14").Value = "  -6.99%  "
$ws.Range("D15").Value = "3.722.64"
$ws.Range("E15").Value = "  -6.18%  "
$ws.Range("E16").Value = "  -7.03%  "
$ws.Range("D17").Value = "3.207.22"
$ws.Range("E17").Value = "  -6.35%  "
$ws.Range("D18").Value = "62.746.01"
$ws.Range("E18").Value = "  -5.05%  "
$ws.Range("D19").Value = "17.09"
$ws.Range("E19").Value = "  -3.87%  "
$ws.Range("D20").Value = "10.90"
$ws.Range("E20").Value = "  -5.14%  "
$ws.Range("D21").Value = "0.950"
$ws.Range("E21").Value = "  -4.96%  "
$ws.Range("D22").Value = "363.76"
$ws.Range("E22").Value = "  -6.03%  "
$ws.Range("E23").Value = "  -3.84%  "
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").Value = "79.81"
$ws.Range("E24").Value = "  -4.92%  "
$ws.Range("B25").Value = "RenderToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D25").Value = "11.00"
$ws.Range("E25").Value = "  -1.87%  "
$ws.Range("D26").Value = "3.87"
$ws.Range("E26").Value = "  +2.48%  "
$ws.Range("D27").Value = "6.03"
$ws.Range("E27").Value = "  -0.68%  "
$ws.Range("E28").Value = "  -5.67%  "
$ws.Range("E29").Value = "  -6.94%  "
$ws.Range("D30").Value = "8.10"
$ws.Range("E30").Value = "  -6.98%  "
$ws.Range("D31").Value = "649.39"
$ws.Range("E31").Value = "  -7.20%  "
$ws.Range("D32").Value = "28.10"
$ws.Range("E33").Value = "  -10.15%  "
$ws.Range("D34").Value = "11.04"
$ws.Range("E34").Value = "  -3.45%  "
$ws.Range("E35").Value = "  -4.60%  "
$ws.Range("D36").Value = "57.80"
$ws.Range("E36").Value = "  -7.22%  "
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("E38").Value = "  -3.12%  "
$ws.Range("E39").Value = "  -4.82%  "
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").Value = "  +0.07%  "
$ws.Range("D41").Value = "0.0₃0689"
$ws.Range("E41").Value = "  +7.41%  "
$ws.Range("D42").Value = "0.122"
$ws.Range("E42").Value = "  -4.90%  "
$ws.Range("D43").Value = "2.850.66"
$ws.Range("E43").Value = "  -3.06%  "
$ws.Range("E44").Value = "  +2.58%  "
$ws.Range("D45").Value = "2.68"
$ws.Range("E45").Value = "  -2.35%  "
$ws.Range("B46").Value = "Stacks"
$ws.Range("C46").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D46").Value = "2.82"
$ws.Range("E46").Value = "  +7.01%  "
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").Value = "0.0386"
$ws.Range("E47").Value = "  -2.67%  "
$ws.Range("E48").Value = "  -10.17%  "
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").Value = "134.99"
$ws.Range("E49").Value = "  +0.43%  "
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").Value = "0.122"
$ws.Range("E50").Value = "  -4.46%  "
$ws.Range("B51").Value = "ApeXProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D51").Value = "2.89"
$ws.Range("E51").Value = "  -1.12%  "
